$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.857.22'
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('D3').Value = '3.114.89'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '171.93'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.16%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '3.111.04'
$ws.Range('E8').Value = '  +1.00%  '
$ws.Range('E9').Value = '  -0.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.44'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.60%  '
$ws.Range('E11').Value = '  -0.47%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.484'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.95'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').Value = '3.631.07'
$ws.Range('E16').Value = '  +1.01%  '
$ws.Range('D17').Value = '66.817.05'
$ws.Range('E17').Value = '  -0.06%  '
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('D19').Value = '3.112.98'
$ws.Range('E19').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.33'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '477.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.03'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('E25').Value = '  +3.33%  '
$ws.Range('E26').Value = '  -3.26%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.87'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.26%  '
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.53'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('E34').Value = '  -7.80%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('E37').Value = '  -2.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '47.19'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '50.20'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E40').Value = '  -3.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.311'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.33%  '
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '2.823.56'
$ws.Range('E44').Value = '  +2.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '382.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('E47').Value = '  -9.26%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '135.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '24.80'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.74%  '
$ws.Range('E51').Value = '  -1.86%  '
